$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of bitcoin-buy data appended after the 2025-05-14 run.
# Column A holds the date as literal text (matching the existing rows'
# "MM/DD/YYYY" inline strings), so format the cell as Text first to stop
# Excel's automatic date recognition from turning it into a date serial.
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "05/14/2025"
$ws.Range("A15").ClearFormats()

$ws.Range("B15").Value = 0.00048211
$ws.Range("C15").Value = 103710.7714007177
$ws.Range("D15").Value = 50
